$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates: volume number and reporting week ---
$ws.Range("A8").Value = "Volume 30   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/12/2023  Through  6/18/2023"

# --- Weekly crime statistics refresh ---
# Row 14
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("C14").Value = 1
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("D14").Value = 1
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 9
$ws.Range("H14").Value = -44.444444444444
$ws.Range("I14").Value = 31
$ws.Range("J14").Value = 32
$ws.Range("K14").Value = -3.125
$ws.Range("L14").Value = -26.190476190476
$ws.Range("M14").Value = -44.642857142857
$ws.Range("N14").Value = -85.024154589372

# Row 15
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 600
$ws.Range("F15").Value = 18
$ws.Range("G15").Value = 13
$ws.Range("H15").Value = 38.461538461538
$ws.Range("I15").Value = 105
$ws.Range("J15").Value = 106
$ws.Range("K15").Value = -0.943396226415
$ws.Range("L15").Value = 11.702127659574
$ws.Range("M15").Value = 5
$ws.Range("N15").Value = -63.286713286713

# Row 16
$ws.Range("C16").Value = 46
$ws.Range("D16").Value = 55
$ws.Range("E16").Value = -16.363636363636
$ws.Range("F16").Value = 187
$ws.Range("G16").Value = 212
$ws.Range("H16").Value = -11.792452830188
$ws.Range("I16").Value = 1086
$ws.Range("J16").Value = 1157
$ws.Range("K16").Value = -6.136560069144
$ws.Range("L16").Value = 23.690205011389
$ws.Range("M16").Value = -28.458498023715
$ws.Range("N16").Value = -85.316387236344

# Row 17
$ws.Range("C17").Value = 100
$ws.Range("D17").Value = 94
$ws.Range("E17").Value = 6.382978723404
$ws.Range("F17").Value = 352
$ws.Range("G17").Value = 353
$ws.Range("H17").Value = -0.28328611898
$ws.Range("I17").Value = 1907
$ws.Range("J17").Value = 1824
$ws.Range("K17").Value = 4.550438596491
$ws.Range("L17").Value = 28.590694538098
$ws.Range("M17").Value = 30.616438356164
$ws.Range("N17").Value = -48.956102783725

# Row 18
$ws.Range("C18").Value = 34
$ws.Range("D18").Value = 52
$ws.Range("E18").Value = -34.615384615384
$ws.Range("F18").Value = 133
$ws.Range("G18").Value = 191
$ws.Range("H18").Value = -30.366492146596
$ws.Range("I18").Value = 919
$ws.Range("J18").Value = 1127
$ws.Range("K18").Value = -18.456078083407
$ws.Range("L18").Value = 8.757396449704
$ws.Range("M18").Value = -28.649068322981
$ws.Range("N18").Value = -82.686510926902

# Row 19
$ws.Range("C19").Value = 112
$ws.Range("D19").Value = 114
$ws.Range("E19").Value = -1.754385964912
$ws.Range("F19").Value = 439
$ws.Range("G19").Value = 467
$ws.Range("H19").Value = -5.995717344753
$ws.Range("I19").Value = 2609
$ws.Range("J19").Value = 2563
$ws.Range("K19").Value = 1.794771751853
$ws.Range("L19").Value = 35.956227201667
$ws.Range("M19").Value = 49.000571102227
$ws.Range("N19").Value = -9.283727399165

# Row 20
$ws.Range("C20").Value = 36
$ws.Range("D20").Value = 33
$ws.Range("E20").Value = 9.090909090909
$ws.Range("F20").Value = 144
$ws.Range("G20").Value = 112
$ws.Range("H20").Value = 28.571428571428
$ws.Range("I20").Value = 790
$ws.Range("J20").Value = 787
$ws.Range("K20").Value = 0.381194409148
$ws.Range("L20").Value = 29.934210526315
$ws.Range("M20").Value = 21.913580246913
$ws.Range("N20").Value = -82.078039927404

# Row 21
$ws.Range("C21").Value = 336
$ws.Range("D21").Value = 350
$ws.Range("E21").Value = -4
$ws.Range("F21").Value = 1278
$ws.Range("G21").Value = 1357
$ws.Range("H21").Value = -5.821665438467
$ws.Range("I21").Value = 7447
$ws.Range("J21").Value = 7596
$ws.Range("K21").Value = -1.961558715113
$ws.Range("L21").Value = 26.887033566195
$ws.Range("M21").Value = 9.177539950153
$ws.Range("N21").Value = -69.248874757401

# Row 22
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = -57.142857142857
$ws.Range("F22").Value = 18
$ws.Range("G22").Value = 23
$ws.Range("H22").Value = -21.739130434782
$ws.Range("I22").Value = 132
$ws.Range("J22").Value = 172
$ws.Range("K22").Value = -23.255813953488
$ws.Range("L22").Value = 14.782608695652
$ws.Range("M22").Value = -32.994923857868
$ws.Range("N22").Value = "***.*"

# Row 23
$ws.Range("C23").Value = 31
$ws.Range("D23").Value = 23
$ws.Range("E23").Value = 34.782608695652
$ws.Range("F23").Value = 135
$ws.Range("G23").Value = 117
$ws.Range("H23").Value = 15.384615384615
$ws.Range("I23").Value = 745
$ws.Range("J23").Value = 677
$ws.Range("K23").Value = 10.044313146233
$ws.Range("L23").Value = 17.693522906793
$ws.Range("M23").Value = 51.422764227642
$ws.Range("N23").Value = "***.*"

# Row 24
$ws.Range("C24").Value = 220
$ws.Range("D24").Value = 262
$ws.Range("E24").Value = -16.030534351145
$ws.Range("F24").Value = 992
$ws.Range("G24").Value = 1045
$ws.Range("H24").Value = -5.071770334928
$ws.Range("I24").Value = 5624
$ws.Range("J24").Value = 5893
$ws.Range("K24").Value = -4.564737824537
$ws.Range("L24").Value = 25.395763656633
$ws.Range("M24").Value = 26.29687850887
$ws.Range("N24").Value = "***.*"

# Row 25
$ws.Range("C25").Value = 132
$ws.Range("D25").Value = 138
$ws.Range("E25").Value = -4.347826086956
$ws.Range("F25").Value = 502
$ws.Range("G25").Value = 580
$ws.Range("H25").Value = -13.448275862069
$ws.Range("I25").Value = 2764
$ws.Range("J25").Value = 2855
$ws.Range("K25").Value = -3.187390542907
$ws.Range("L25").Value = 40.375825292026
$ws.Range("M25").Value = -22.836404243439
$ws.Range("N25").Value = "***.*"

# Row 26
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 175
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 156
$ws.Range("J26").Value = 171
$ws.Range("K26").Value = -8.771929824561
$ws.Range("L26").Value = -6.586826347305
$ws.Range("M26").Value = "***.*"
$ws.Range("N26").Value = "***.*"

# Row 27
$ws.Range("C27").Value = 18
$ws.Range("D27").Value = 15
$ws.Range("E27").Value = 20
$ws.Range("F27").Value = 55
$ws.Range("G27").Value = 53
$ws.Range("H27").Value = 3.77358490566
$ws.Range("I27").Value = 286
$ws.Range("J27").Value = 279
$ws.Range("K27").Value = 2.508960573476
$ws.Range("L27").Value = -8.038585209003
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"

# Row 28
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 20
$ws.Range("F28").Value = 19
$ws.Range("G28").Value = 35
$ws.Range("H28").Value = -45.714285714285
$ws.Range("I28").Value = 106
$ws.Range("J28").Value = 134
$ws.Range("K28").Value = -20.895522388059
$ws.Range("L28").Value = -36.144578313253
$ws.Range("M28").Value = -47
$ws.Range("N28").Value = -87.745664739884

# Row 29
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = 25
$ws.Range("F29").Value = 16
$ws.Range("G29").Value = 28
$ws.Range("H29").Value = -42.857142857142
$ws.Range("I29").Value = 91
$ws.Range("J29").Value = 114
$ws.Range("K29").Value = -20.175438596491
$ws.Range("L29").Value = -36.363636363636
$ws.Range("M29").Value = -43.125
$ws.Range("N29").Value = -88.392857142857

# Row 30
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -50
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = -80
$ws.Range("I30").Value = 30
$ws.Range("J30").Value = 37
$ws.Range("K30").Value = -18.918918918918
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = "***.*"
$ws.Range("N30").Value = "***.*"
